# Generate Report for Handoff
# Adds a new handed-off file (a7f5cb4a-ba4e-4ade-83c9-de6fc6905ff4) as a new
# row on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$guid = "a7f5cb4a-ba4e-4ade-83c9-de6fc6905ff4"
$mdDisplay = "$guid.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/db7738a60eca59c66d9d14e1ff1c962483df9611/e2e/$mdDisplay"

# ---------------------------------------------------------------------------
# Overview sheet - new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdDisplay
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-28-20 10:28:40"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $mdDisplay)
$wsOverview.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# zh-cn sheet - new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfDisplay = "$guid.a555739c631a8e7de1181c1aa407048433d7e8c9.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01ba3871b64bc2423b6e20a4868d0f8de26ef42c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfDisplay"

$wsZhCn.Range("A3").Value = $mdDisplay
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = $zhXlfDisplay
$wsZhCn.Range("E3").Value = "2016-03-20 10:28:37"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, "", "", $mdDisplay)
$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $mdUrl, "", "", ".md")
$wsZhCn.Range("B3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhXlfUrl, "", "", $zhXlfDisplay)
$wsZhCn.Range("D3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# de-de sheet - new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfDisplay = "$guid.a555739c631a8e7de1181c1aa407048433d7e8c9.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/faca294e54ef6ea777d1f2c2822ed2e9635ad55c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfDisplay"

$wsDeDe.Range("A3").Value = $mdDisplay
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = $deXlfDisplay
$wsDeDe.Range("E3").Value = "2016-03-20 10:28:40"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, "", "", $mdDisplay)
$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDeDe.Range("B3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deXlfUrl, "", "", $deXlfDisplay)
$wsDeDe.Range("D3").Style = "HyperLink"

Write-Host "Report generated for handoff of $mdDisplay"
